# Flip the page from portrait to landscape (wdOrientLandscape = 1).
# Word's PageSetup.Orientation setter swaps PageWidth/PageHeight for us,
# matching the pgSz w/h swap + orient="landscape" seen in the target XML.
$d = $word.ActiveDocument
$d.PageSetup.Orientation = 1
